# Inventory Management System final code review.pptx
# "striped out screen shots of web pages"
#
# 1) Remove the 7 web-page screenshot slides (originally slides 4-10,
#    i.e. the picture-only slides sitting between "Interface" and
#    "Storage"). Deleting slide #4 repeatedly shifts the remaining
#    slides up into place each time.
# 2) On the "Storage" slide (now slide #4), fix the "Website
#    successfully connects to the server when hosted localy" line to
#    read "...hosted locally, But is having trouble with editing ".

$p = $ppt.ActivePresentation

# --- 1) Delete the 7 consecutive screenshot slides --------------------
for ($i = 0; $i -lt 7; $i++) {
    $p.Slides.Item(4).Delete()
}

# --- 2) Fix the wording on the Storage slide's "Final State" box ------
$storageSlide = $p.Slides.Item(4)
$finalStateBox = $storageSlide.Shapes.Item(5)
$tr = $finalStateBox.TextFrame.TextRange
$para = $tr.Paragraphs(3, 1)

$run1 = $para.Runs(1, 1)
$run1.Text = "Website successfully connects to the server when hosted locally, But is having trouble with editing "

$run2 = $para.Runs(2, 1)
$run2.Text = ""
